# Updated cryptos list on Thu May 18 19:41:19 UTC 2023 with GitHub Actions
#
# This script refreshes the "Price" (column D) and "Volume(1h)" (column E)
# figures for every coin row on the sheet, and also fixes up three rows
# whose coin/link (columns B and C) had been reordered in the source feed
# (EthereumClassic <-> LidoDAOToken around rows 27-28, and
#  MXToken / FraxShare / RenderToken around rows 40-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry below corresponds to one data row (row 1 is the header).
# B/C are only present where the coin name/link actually changed;
# D (Price) and E (Volume(1h)) are refreshed for (almost) every row.
$updates = @(
    @{ Row=2;  D='26.709.86';    E='  -2.29%  ' }
    @{ Row=3;  D='1.795.16';     E='  -1.76%  ' }
    @{ Row=4;  D='0.9995';       E='  +0.00%  ' }
    @{ Row=5;  D='308.08';       E='  -1.74%  ' }
    @{ Row=6;  D='0.9988';       E='  -0.04%  ' }
    @{ Row=7;  D='0.4561';       E='  +1.66%  ' }
    @{ Row=8;  D='0.3709';       E='  -1.86%  ' }
    @{ Row=9;  D='0.07231';      E='  -3.75%  ' }
    @{ Row=10; D='0.8553';       E='  -4.37%  ' }
    @{ Row=11; D='20.39';        E='  -3.26%  ' }
    @{ Row=12; D='1.802.51';     E='  -1.34%  ' }
    @{ Row=13; D='5.299';        E='  -1.97%  ' }
    @{ Row=14; D='6.494';        E='  -4.19%  ' }
    @{ Row=15; D='0.07033';      E='  -1.29%  ' }
    @{ Row=16; D='90.16';        E='  -4.58%  ' }
    @{ Row=17; D='0.9995';       E='  -0.04%  ' }
    @{ Row=18; D='0.000008632';  E='  -2.26%  ' }
    @{ Row=19; D='0.9990';       E='  -0.02%  ' }
    @{ Row=20; D='14.60' }
    @{ Row=21; D='26.733.22';    E='  -2.29%  ' }
    @{ Row=22; D='5.288';        E='  +0.04%  ' }
    @{ Row=23; D='10.59';        E='  -3.52%  ' }
    @{ Row=24; D='2.018.13';     E='  -1.56%  ' }
    @{ Row=25; D='1.907';        E='  -4.25%  ' }
    @{ Row=26; D='149.76';       E='  -1.43%  ' }
    @{ Row=27; B='LidoDAOToken';    C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo';           D='2.153'; E='  -13.02%  ' }
    @{ Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc';    D='18.15'; E='  -2.66%  '  }
    @{ Row=29; D='5.197';        E='  -3.51%  ' }
    @{ Row=30; D='113.98';       E='  -3.74%  ' }
    @{ Row=31; D='0.08833';      E='  -0.20%  ' }
    @{ Row=32; D='0.7576';       E='  -2.76%  ' }
    @{ Row=33; D='1.157';        E='  -3.37%  ' }
    @{ Row=34; D='4.443';        E='  -3.02%  ' }
    @{ Row=35; D='2.883';        E='  -0.06%  ' }
    @{ Row=36; D='0.9982';       E='  -0.01%  ' }
    @{ Row=37; D='1.109';        E='  -0.44%  ' }
    @{ Row=38; D='0.01939';      E='  -2.77%  ' }
    @{ Row=39; D='0.05200';      E='  -2.38%  ' }
    @{ Row=40; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr';  D='2.369'; E='  +4.39%  ' }
    @{ Row=41; B='MXToken';     C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx';          D='2.884'; E='  +0.82%  ' }
    @{ Row=42; B='FraxShare';   C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs';          D='7.117'; E='  -4.24%  ' }
    @{ Row=43; D='0.5227';       E='  -2.61%  ' }
    @{ Row=44; D='0.1643';       E='  -5.45%  ' }
    @{ Row=45; D='8.480';        E='  -3.96%  ' }
    @{ Row=46; D='0.5000';       E='  -2.99%  ' }
    @{ Row=47; D='10.25';        E='  -5.02%  ' }
    @{ Row=48; D='104.24';       E='  -2.59%  ' }
    @{ Row=49; D='0.9977';       E='  -0.06%  ' }
    @{ Row=50; D='1.646';        E='  -3.70%  ' }
    @{ Row=51; D='0.06278';      E='  -1.66%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Range("C$row").Value = $u.C
    }

    if ($u.ContainsKey('D')) {
        $priceCell = $ws.Range("D$row")
        # The Price column holds plain text (e.g. "26.709.86", "0.9995").
        # Some of the new values parse as ordinary numbers, which would
        # make Excel silently convert them to numeric cells. Forcing a
        # text number format before assigning keeps them as text, just
        # like the original (and like the sibling D values that contain
        # a second "." and so can never be mistaken for numbers).
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.D
    }

    if ($u.ContainsKey('E')) {
        $ws.Range("E$row").Value = $u.E
    }
}
